$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting the existing rows 4-86 down to 5-87.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new daily record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R stay the same as the surrounding rows;
# only D, J, K, M, P receive fresh values (L keeps the same value it had before).
$ws.Range("A4").Value2 = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value2 = 45245
$ws.Range("E4").Value2 = 10
$ws.Range("F4").Value2 = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value2 = 100
$ws.Range("K4").Value2 = 2000
$ws.Range("L4").Value2 = 2000
$ws.Range("M4").Value2 = 2000
$ws.Range("N4").Value = "$/kilo"
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value2 = 2000
$ws.Range("Q4").Value2 = 1
$ws.Range("R4").Value = "Hortaliza"
